# Bugfix for missing municipalities
# 1) Grande-Anse (2005 tax base data) was incorrectly duplicated from the
#    Grand Manan row; restore its real values and slightly-adjusted G553.
# 2) Four municipalities (Bathurst, Campbellton, Dieppe, Edmundston) were
#    incorrectly excluded from the 2020 data; insert them back in at the
#    top of the 2020 block, shifting the rest of the 2020 rows down by 4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix 1: 2005 Grande-Anse / Grand Manan rows -----------------------
$ws.Cells.Item(553, 7).Value = 51.98573333333333

$ws.Cells.Item(554, 2).Value  = "Grande-Anse"
$ws.Cells.Item(554, 3).Value  = 1.341
$ws.Cells.Item(554, 4).Value  = 0.09100000000000001
$ws.Cells.Item(554, 5).Value  = 0.8130480656506448
$ws.Cells.Item(554, 6).Value  = -0.688926143024619
$ws.Cells.Item(554, 7).Value  = 36.42482415005862
$ws.Cells.Item(554, 10).Value = 853

# --- Fix 2: insert the 4 missing 2020 municipalities -------------------
$ws.Rows("2018:2021").Insert()

# Row 2018: Bathurst
$ws.Cells.Item(2018, 1).Value  = 2020
$ws.Cells.Item(2018, 2).Value  = "Bathurst"
$ws.Cells.Item(2018, 3).Value  = 1.775
$ws.Cells.Item(2018, 4).Value  = 0.3607543918634951
$ws.Cells.Item(2018, 5).Value  = 1.818663276456249
$ws.Cells.Item(2018, 6).Value  = 0.5694793645456836
$ws.Cells.Item(2018, 7).Value  = 90.7007474993696
$ws.Cells.Item(2018, 8).Value  = $false
$ws.Cells.Item(2018, 9).Value  = $false
$ws.Cells.Item(2018, 10).Value = 11897

# Row 2019: Campbellton
$ws.Cells.Item(2019, 1).Value  = 2020
$ws.Cells.Item(2019, 2).Value  = "Campbellton"
$ws.Cells.Item(2019, 3).Value  = 1.7763
$ws.Cells.Item(2019, 4).Value  = 0.3826215313090222
$ws.Cells.Item(2019, 5).Value  = 1.774466947551939
$ws.Cells.Item(2019, 6).Value  = 0.6045397355804153
$ws.Cells.Item(2019, 7).Value  = 87.40352317303503
$ws.Cells.Item(2019, 8).Value  = $false
$ws.Cells.Item(2019, 9).Value  = $true
$ws.Cells.Item(2019, 10).Value = 6883

# Row 2020: Dieppe
$ws.Cells.Item(2020, 1).Value  = 2020
$ws.Cells.Item(2020, 2).Value  = "Dieppe"
$ws.Cells.Item(2020, 3).Value  = 1.6295
$ws.Cells.Item(2020, 4).Value  = 0.2344724235739048
$ws.Cells.Item(2020, 5).Value  = 2.087807004412228
$ws.Cells.Item(2020, 6).Value  = 0.2216123148439962
$ws.Cells.Item(2020, 7).Value  = 128.9160180822566
$ws.Cells.Item(2020, 8).Value  = $false
$ws.Cells.Item(2020, 9).Value  = $true
$ws.Cells.Item(2020, 10).Value = 25384

# Row 2021: Edmundston
$ws.Cells.Item(2021, 1).Value  = 2020
$ws.Cells.Item(2021, 2).Value  = "Edmundston"
$ws.Cells.Item(2021, 3).Value  = 1.635
$ws.Cells.Item(2021, 4).Value  = 0.33544794933655
$ws.Cells.Item(2021, 5).Value  = 1.806330398069964
$ws.Cells.Item(2021, 6).Value  = 0.7111211097708079
$ws.Cells.Item(2021, 7).Value  = 87.5019722557298
$ws.Cells.Item(2021, 8).Value  = $false
$ws.Cells.Item(2021, 9).Value  = $false
$ws.Cells.Item(2021, 10).Value = 16580

# --- Grow the "Frame0" table / autofilter to cover the new rows --------
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:J2110"))
